# amendStatusOfNonMemberViaDueDate feature:
#  - rename the old "sqlCount" sheet to "amendStatusOfNonMemberViaDueDate"
#    (Excel COM enforces the 31-char sheet-name limit, so the name is
#    truncated to fit) and populate it with the emailId test data.
#  - add a fresh "sqlCount" sheet after it, restoring the original
#    sqlCount content.
#  - add a row to createRegionalEvent with the test event name.
#  - make addBrandNewVisitor the active/selected sheet.

$wb = $excel.ActiveWorkbook

# 1. Add a row to createRegionalEvent (eventName header already present in A1)
$regionalEventSheet = $wb.Worksheets.Item("createRegionalEvent")
$regionalEventSheet.Range("A2").Value = "TestAutomation20191211060135"

# 2. Rename the existing sqlCount sheet to the new feature sheet name and
#    replace its content with the emailId test data.
$sqlCountSheet = $wb.Worksheets.Item("sqlCount")
$sqlCountSheet.Name = "amendStatusOfNonMemberViaDueDat"
$sqlCountSheet.Range("A1").Value = "emailId"
$sqlCountSheet.Range("A2").Value = "shanthibni+32@gmail.com"

# 3. Add a brand new sqlCount sheet right after it, with the original
#    sqlCount content. (The counts are stored as text in the original
#    workbook, so force a text number-format before assigning them.)
$newSqlCountSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sqlCountSheet)
$newSqlCountSheet.Name = "sqlCount"
$newSqlCountSheet.Range("A1").Value = "sqlRecordCount"
$newSqlCountSheet.Range("B1").Value = "sqlColCount"
$newSqlCountSheet.Range("A2:B2").NumberFormat = "@"
$newSqlCountSheet.Range("A2").Value = "252"
$newSqlCountSheet.Range("B2").Value = "1"

# 4. Make addBrandNewVisitor the active/selected sheet.
$firstSheet = $wb.Worksheets.Item("addBrandNewVisitor")
$firstSheet.Activate()
